# Update weekly Madigan bike-share ridership "Riders" (C) and "Average" (D)
# values with the newly computed hours for the 2017 week 19 report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Row 2 - Monday 08 May 2017
$ws.Range("C2").Value = 187
$ws.Range("D2").Value = 226.12

# Row 3 - Tuesday 09 May 2017
$ws.Range("C3").Value = 244
$ws.Range("D3").Value = 220.96

# Row 4 - Wednesday 10 May 2017
$ws.Range("C4").Value = 207
$ws.Range("D4").Value = 213.68

# Row 5 - Thursday 11 May 2017
$ws.Range("C5").Value = 218
$ws.Range("D5").Value = 234.56

# Row 6 - Friday 12 May 2017
$ws.Range("C6").Value = 191
$ws.Range("D6").Value = 239.61

# Row 7 - Saturday 13 May 2017
$ws.Range("C7").Value = 93
$ws.Range("D7").Value = 111.45

# Row 8 - Sunday 14 May 2017
$ws.Range("C8").Value = 105
$ws.Range("D8").Value = 91.19
